# LegacyCollectionsTemplate.docx fix-up
#
# 1. "Captain: {{" + "captain-name" + "}}" (3 runs) collapses into a single
#    run "Captain: {{captain-name}}", and a brand-new paragraph is inserted
#    right after it for the (until now missing) first-officer merge field,
#    with the placeholder's "c" wrapped in a DDE-link bookmark (an artifact
#    of how this text was pasted into Word).
# 2. The table's {{officer-name}} / {{officer-rank}} / {{officer-uniform}}
#    placeholders (each split across 3 runs: "{{", name, "}}") collapse to
#    one run apiece.
# 3. {{service-ship-name}} (again split across 3 runs) collapses to one run.

$d = $word.ActiveDocument

function Get-VisibleLength($text) {
    # Table-cell ranges (and some paragraph ranges) report a trailing
    # paragraph/cell mark in .Text ("\r", or "\r\a" for an end-of-cell
    # marker) that isn't part of the real Start/End addressable span we
    # want to touch. Strip trailing control characters to get the true
    # visible-text length.
    $len = $text.Length
    while ($len -gt 0 -and [int][char]$text.Substring($len - 1, 1) -lt 32) {
        $len = $len - 1
    }
    return $len
}

function Merge-BraceRuns($range) {
    # $range's visible text is exactly "{{<name>}}" split across 3 runs
    # ("{{", "<name>", "}}"). Delete everything after the opening "{{" and
    # re-append it to the first run so the first run's rPr (formatting)
    # survives the merge, instead of Find/Replace's whole-match clobber.
    $visibleLen = Get-VisibleLength $range.Text
    $start = $range.Start
    $headEnd = $start + 2
    $tailEnd = $start + $visibleLen
    $tailText = $d.Range($headEnd, $tailEnd).Text
    $d.Range($headEnd, $tailEnd).Delete()
    $d.Range($start, $headEnd).InsertAfter($tailText)
}

# ---------------------------------------------------------------------
# 1. Captain paragraph: merge its 3 runs, then split off a new paragraph
#    carrying "{{captain-first-officer-name}}" with a bookmark around the
#    lone "c" right after the opening braces.
# ---------------------------------------------------------------------
$captainPara = $d.Paragraphs(1)
Merge-BraceRuns $captainPara.Range

# Re-fetch the paragraph (its Range may be stale after the edit) and
# insert a brand-new paragraph right after it, matching its style.
$captainPara = $d.Paragraphs(1)
$captainPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs(2)
$newPara.Range.InsertAfter("{{captain-first-officer-name}}")

# Wrap just the "c" (right after the opening "{{") in a bookmark, matching
# the stray __DdeLink bookmark left behind by whatever pasted this text in.
$newParaStart = $newPara.Range.Start
$cRange = $d.Range($newParaStart + 2, $newParaStart + 3)
$d.Bookmarks.Add("__DdeLink__44_427613671", $cRange) | Out-Null

# ---------------------------------------------------------------------
# 2. Collapse the table's "{{" / name / "}}" 3-run placeholders to one
#    run each.
# ---------------------------------------------------------------------
$table = $d.Tables(1)
Merge-BraceRuns $table.Cell(2, 1).Range
Merge-BraceRuns $table.Cell(2, 2).Range
Merge-BraceRuns $table.Cell(2, 3).Range

# ---------------------------------------------------------------------
# 3. Same collapse for the {{service-ship-name}} list item.
# ---------------------------------------------------------------------
$searchRange = $d.Content
$found = $searchRange.Find.Execute("{{service-ship-name}}")
if ($found) {
    Merge-BraceRuns $searchRange
}
